$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff.
# D-column (Price) values are forced to Text first: many of the new prices
# (e.g. "234.67") would otherwise be auto-parsed by Excel as numbers, which
# does not match the workbook's original inlineStr/text cell type.
$priceCells = @{
    "D2" = "37.373.49"
    "D3" = "2.064.85"
    "D5" = "234.67"
    "D7" = "57.95"
    "D9" = "0.381"
    "D10" = "58.87"
    "D13" = "2.370.53"
    "D14" = "14.52"
    "D15" = "21.10"
    "D17" = "5.19"
    "D18" = "2.079.71"
    "D19" = "37.582.67"
    "D21" = "70.39"
    "D23" = "227.13"
    "D25" = "2.45"
    "D27" = "165.19"
    "D29" = "8.88"
    "D30" = "19.20"
    "D32" = "0.118"
    "D33" = "4.52"
    "D38" = "3.37"
    "D41" = "2.96"
    "D42" = "0.0968"
    "D43" = "4.41"
    "D44" = "96.11"
    "D45" = "1.453.62"
    "D48" = "15.79"
}
foreach ($ref in $priceCells.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $priceCells[$ref]
    $rng.Style = "Normal"
}

# Remaining text cells (coin name / link / volume %) already round-trip as text.
$textCells = @{
    "E2" = "  +2.24%  "
    "E4" = "  -0.06%  "
    "E5" = "  +0.76%  "
    "E6" = "  +2.27%  "
    "E7" = "  +5.34%  "
    "E8" = "  -0.03%  "
    "E9" = "  +2.85%  "
    "E10" = "  +1.34%  "
    "E11" = "  +1.91%  "
    "E12" = "  +2.84%  "
    "E13" = "  +3.49%  "
    "E14" = "  +2.26%  "
    "E15" = "  +3.95%  "
    "E16" = "  +2.38%  "
    "E17" = "  +2.27%  "
    "E18" = "  +4.17%  "
    "E19" = "  +3.04%  "
    "E20" = "  +16.25%  "
    "E21" = "  +3.92%  "
    "E22" = "  +1.25%  "
    "E23" = "  +2.29%  "
    "E24" = "  -0.08%  "
    "E25" = "  +2.18%  "
    "E26" = "  +0.81%  "
    "E27" = "  +1.96%  "
    "E28" = "  +12.70%  "
    "E29" = "  +2.43%  "
    "E30" = "  +1.80%  "
    "E31" = "  +0.66%  "
    "E32" = "  +0.99%  "
    "E33" = "  +3.15%  "
    "E34" = "  +2.67%  "
    "E35" = "  +8.94%  "
    "E36" = "  +6.74%  "
    "E37" = "  +0.00%  "
    "E38" = "  +0.34%  "
    "E39" = "  +0.82%  "
    "E40" = "  +3.95%  "
    "B41" = "HuobiToken"
    "C41" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "E41" = "  -1.55%  "
    "B42" = "Cronos"
    "C42" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "E42" = "  +2.96%  "
    "E43" = "  +21.53%  "
    "E44" = "  +7.80%  "
    "E45" = "  -0.24%  "
    "E46" = "  +3.92%  "
    "E47" = "  +5.52%  "
    "E48" = "  +3.98%  "
    "E49" = "  +4.05%  "
    "E50" = "  +6.57%  "
    "E51" = "  +2.04%  "
}
foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).Value = $textCells[$ref]
}
